$wb = $excel.ActiveWorkbook

# --- Insert a new "Todo" worksheet as the first sheet in the workbook ---
$firstSheet = $wb.Worksheets.Item(1)
$todo = $wb.Worksheets.Add($firstSheet)
$todo.Name = "Todo"

# --- Header row ---
$todo.Range("A1").Value = "Task #"
$todo.Range("B1").Value = "Description"
$todo.Range("C1").Value = "Owner"
$todo.Range("D1").Value = "Status"

# --- Task rows 2-5 (entered in natural order) ---
$todo.Range("A2").Value = 1
$todo.Range("B2").Value = "Define typeDef to establish database"
$todo.Range("A3").Value = 2
$todo.Range("B3").Value = "Build app component"
$todo.Range("A4").Value = 3
$todo.Range("B4").Value = "Build navigation component"
$todo.Range("A5").Value = 4
$todo.Range("B5").Value = "build login modal"

# Rows 7-9 are populated before row 6 so the shared-string table ends up
# in the same order as the authored workbook (build logout modal was
# typed after these three, then the row was moved above them).
$todo.Range("A7").Value = 6
$todo.Range("B7").Value = "build image sharing component"
$todo.Range("A8").Value = 7
$todo.Range("B8").Value = "build image creation component"
$todo.Range("A9").Value = 8
$todo.Range("B9").Value = "build leaderboard component"

$todo.Range("A6").Value = 5
$todo.Range("B6").Value = "build logout modal"

$todo.Range("A10").Value = 9
$todo.Range("B10").Value = "create addSharedImage mutation"
$todo.Range("A11").Value = 10
$todo.Range("B11").Value = "create removeSharedImage mutation"
$todo.Range("A12").Value = 11
$todo.Range("B12").Value = "create addCreatedImage mutation"
$todo.Range("A13").Value = 12
$todo.Range("B13").Value = "create removeCreatedImage mutation"
$todo.Range("A14").Value = 13
$todo.Range("B14").Value = "create registerUser mutation"

# --- Remaining numbered rows (no description yet) ---
$todo.Range("A15").Value = 14
$todo.Range("A16").Value = 15
$todo.Range("A17").Value = 16
$todo.Range("A18").Value = 17
$todo.Range("A19").Value = 18
$todo.Range("A20").Value = 19
$todo.Range("A21").Value = 20
$todo.Range("A22").Value = 21
$todo.Range("A23").Value = 22
$todo.Range("A24").Value = 23
$todo.Range("A25").Value = 24
$todo.Range("A26").Value = 25
$todo.Range("A27").Value = 26
$todo.Range("A28").Value = 27
$todo.Range("A29").Value = 28
$todo.Range("A30").Value = 29
$todo.Range("A31").Value = 30

# --- Column B is widened to fit the task descriptions ---
$todo.Columns.Item(2).ColumnWidth = 115.35

# --- Minor selection changes the author made on a couple of the other
# sheets while working in this session ---
$appComponent = $wb.Worksheets.Item("App Component")
$appComponent.Activate()
$appComponent.Range("F37").Select() | Out-Null

$loginModal = $wb.Worksheets.Item("Login Modal")
$loginModal.Activate()
$loginModal.Range("D44").Select() | Out-Null

# --- Leave the new Todo sheet active/selected, as in the saved file ---
$todo.Activate()
$todo.Range("B4:B5").Select() | Out-Null
